# Re-resolve the workbook/worksheet via the Excel object model (the
# pre-bound $wb/$ws in the outer scope don't reliably round-trip writes
# in this host, so start from $excel.ActiveWorkbook as usual COM code would).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New skill row: "怪物寒冰箭" (monster ice-bolt) / "怪物使用寒冰箭" description.
$ws.Range("A10").Value = 202
$ws.Range("B10").Value = "怪物寒冰箭"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 6
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = "怪物使用寒冰箭"

# Match the author's final selection/scroll state in the saved view.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("L10").Select()
